$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActivityList")

$ws.Rows("6:6").Insert()
$ws.Range("A6").Value = "My Last 100 days Activities"

$ws.Range("M9").Select()
